# Auto-generated edit script: update profit/price figures across all sheets
# per the scheduled runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 78924.64
$ws.Range("J9").Value = 532.8333
$ws.Range("L9").Value = 532.8333
$ws.Range("N9").Value = -870.8333
$ws.Range("H15").Value = 1063.3846
$ws.Range("I15").Value = 1063.3846
$ws.Range("K15").Value = 3190.1538
$ws.Range("M15").Value = -3021.1538
$ws.Range("H17").Value = 1818.4762
$ws.Range("J17").Value = 1818.4762
$ws.Range("L17").Value = 5455.4286
$ws.Range("N17").Value = -5791.4286
$ws.Range("H28").Value = 2892.4546
$ws.Range("I28").Value = 2270.8215
$ws.Range("J28").Value = 6373.6
$ws.Range("K28").Value = 2270.8215
$ws.Range("L28").Value = 6373.6
$ws.Range("M28").Value = -1785.8215
$ws.Range("N28").Value = -7343.6
$ws.Range("H69").Value = 12710.4
$ws.Range("I69").Value = 12768.5
$ws.Range("J69").Value = 12671.667
$ws.Range("K69").Value = 38305.5
$ws.Range("L69").Value = 38015.001
$ws.Range("M69").Value = -37431.5
$ws.Range("N69").Value = -39763.001
$ws.Range("H72").Value = 12710.4
$ws.Range("I72").Value = 12768.5
$ws.Range("J72").Value = 12671.667
$ws.Range("K72").Value = 114916.5
$ws.Range("L72").Value = 114045.003
$ws.Range("M72").Value = -110548.5
$ws.Range("N72").Value = -122781.003
$ws.Range("H86").Value = 3187.0715
$ws.Range("I86").Value = 2692.6365
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2692.6365
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1569.6365
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 3187.0715
$ws.Range("I89").Value = 2692.6365
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 13463.1825
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -7847.182500000001
$ws.Range("N89").Value = -36232
$ws.Range("H103").Value = 1024.3334
$ws.Range("I103").Value = 767
$ws.Range("K103").Value = 2301
$ws.Range("M103").Value = -1715
$ws.Range("H125").Value = 2452
$ws.Range("I125").Value = 1174.5555
$ws.Range("J125").Value = 4751.4
$ws.Range("K125").Value = 10570.9995
$ws.Range("L125").Value = 42762.6
$ws.Range("M125").Value = -8110.9995
$ws.Range("N125").Value = -47682.6
$ws.Range("H129").Value = 1067.3334
$ws.Range("I129").Value = 812
$ws.Range("K129").Value = 2436
$ws.Range("M129").Value = 2564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 740.8
$ws.Range("I2").Value = 600.96155
$ws.Range("K2").Value = 600.96155
$ws.Range("M2").Value = -487.96155
$ws.Range("H61").Value = 10511.546
$ws.Range("I61").Value = 3955.5557
$ws.Range("J61").Value = 40013.5
$ws.Range("K61").Value = 3955.5557
$ws.Range("L61").Value = 40013.5
$ws.Range("M61").Value = -3743.5557
$ws.Range("N61").Value = -40437.5
$ws.Range("H74").Value = 2143.75
$ws.Range("I74").Value = 1465.125
$ws.Range("K74").Value = 1465.125
$ws.Range("M74").Value = -591.125
$ws.Range("H77").Value = 2143.75
$ws.Range("I77").Value = 1465.125
$ws.Range("K77").Value = 7325.625
$ws.Range("M77").Value = -2957.625
$ws.Range("H116").Value = 740.8
$ws.Range("I116").Value = 600.96155
$ws.Range("K116").Value = 600.96155
$ws.Range("M116").Value = 1693.03845
$ws.Range("H132").Value = 3651.4075
$ws.Range("I132").Value = 3464.3696
$ws.Range("J132").Value = 4726.875
$ws.Range("K132").Value = 10393.1088
$ws.Range("L132").Value = 14180.625
$ws.Range("M132").Value = -7863.1088
$ws.Range("N132").Value = -19240.625
$ws.Range("H133").Value = 94998.28
$ws.Range("J133").Value = 94998.28
$ws.Range("L133").Value = 94998.28
$ws.Range("N133").Value = -100058.28
$ws.Range("H136").Value = 10511.546
$ws.Range("I136").Value = 3955.5557
$ws.Range("J136").Value = 40013.5
$ws.Range("K136").Value = 11866.6671
$ws.Range("L136").Value = 120040.5
$ws.Range("M136").Value = -9316.667099999999
$ws.Range("N136").Value = -125140.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 740.8
$ws.Range("I3").Value = 600.96155
$ws.Range("K3").Value = 600.96155
$ws.Range("M3").Value = -486.96155
$ws.Range("H86").Value = 5355.2334
$ws.Range("I86").Value = 5328.7
$ws.Range("J86").Value = 5408.3
$ws.Range("K86").Value = 5328.7
$ws.Range("L86").Value = 5408.3
$ws.Range("M86").Value = -4205.7
$ws.Range("N86").Value = -7654.3
$ws.Range("H89").Value = 5355.2334
$ws.Range("I89").Value = 5328.7
$ws.Range("J89").Value = 5408.3
$ws.Range("K89").Value = 26643.5
$ws.Range("L89").Value = 27041.5
$ws.Range("M89").Value = -21027.5
$ws.Range("N89").Value = -38273.5
$ws.Range("H105").Value = 4614
$ws.Range("I105").Value = 4269.25
$ws.Range("J105").Value = 4889.8
$ws.Range("K105").Value = 4269.25
$ws.Range("L105").Value = 4889.8
$ws.Range("M105").Value = -2522.25
$ws.Range("N105").Value = -8383.799999999999
$ws.Range("H134").Value = 9536.897000000001
$ws.Range("I134").Value = 9294.556
$ws.Range("J134").Value = 10082.167
$ws.Range("K134").Value = 27883.668
$ws.Range("L134").Value = 30246.501
$ws.Range("M134").Value = -25348.668
$ws.Range("N134").Value = -35316.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9332.333000000001
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 9332.333000000001
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H132").Value = 2794.423
$ws.Range("I132").Value = 2409.5715
$ws.Range("K132").Value = 7228.7145
$ws.Range("M132").Value = -4698.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 565.6667
$ws.Range("I5").Value = 565.6667
$ws.Range("K5").Value = 1697.0001
$ws.Range("M5").Value = -1585.0001
$ws.Range("H107").Value = 988.8333
$ws.Range("J107").Value = 1019.63635
$ws.Range("L107").Value = 3058.90905
$ws.Range("N107").Value = -6898.90905
$ws.Range("H111").Value = 6129.6665
$ws.Range("I111").Value = 2400
$ws.Range("K111").Value = 7200
$ws.Range("M111").Value = -4133
$ws.Range("H115").Value = 1732.25
$ws.Range("I115").Value = 646.3333
$ws.Range("K115").Value = 1938.9999
$ws.Range("M115").Value = -763.9999
$ws.Range("H131").Value = 27779268
$ws.Range("I131").Value = 18519504
$ws.Range("J131").Value = 33335126
$ws.Range("K131").Value = 55558512
$ws.Range("L131").Value = 100005378
$ws.Range("M131").Value = -55553472
$ws.Range("N131").Value = -100015458
$ws.Range("H132").Value = 6061.8
$ws.Range("I132").Value = 5439.6665
$ws.Range("J132").Value = 6995
$ws.Range("K132").Value = 48956.9985
$ws.Range("L132").Value = 62955
$ws.Range("M132").Value = -46426.9985
$ws.Range("N132").Value = -68015
$ws.Range("H135").Value = 565.6667
$ws.Range("I135").Value = 565.6667
$ws.Range("K135").Value = 5091.0003
$ws.Range("M135").Value = -2556.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9981.933999999999
$ws.Range("I70").Value = 10898.2
$ws.Range("K70").Value = 10898.2
$ws.Range("M70").Value = -10628.2
$ws.Range("H73").Value = 9981.933999999999
$ws.Range("I73").Value = 10898.2
$ws.Range("K73").Value = 10898.2
$ws.Range("M73").Value = -9962.200000000001
$ws.Range("H80").Value = 3517.9092
$ws.Range("I80").Value = 3232.2856
$ws.Range("K80").Value = 3232.2856
$ws.Range("M80").Value = -2234.2856
$ws.Range("H83").Value = 3517.9092
$ws.Range("I83").Value = 3232.2856
$ws.Range("K83").Value = 16161.428
$ws.Range("M83").Value = -11169.428
$ws.Range("H107").Value = 946.3
$ws.Range("I107").Value = 946.3
$ws.Range("K107").Value = 946.3
$ws.Range("M107").Value = 973.7
$ws.Range("H122").Value = 5342.1377
$ws.Range("I122").Value = 4843.3887
$ws.Range("J122").Value = 6158.273
$ws.Range("K122").Value = 14530.1661
$ws.Range("L122").Value = 18474.819
$ws.Range("M122").Value = -12080.1661
$ws.Range("N122").Value = -23374.819
$ws.Range("H132").Value = 3232.5173
$ws.Range("J132").Value = 5499.7144
$ws.Range("L132").Value = 16499.1432
$ws.Range("N132").Value = -21559.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3607.1628
$ws.Range("I132").Value = 3252.5386
$ws.Range("J132").Value = 4149.5293
$ws.Range("K132").Value = 9757.6158
$ws.Range("L132").Value = 12448.5879
$ws.Range("M132").Value = -7227.6158
$ws.Range("N132").Value = -17508.5879
$ws.Range("H136").Value = 3009.7368
$ws.Range("I136").Value = 2716.087
$ws.Range("K136").Value = 8148.261
$ws.Range("M136").Value = -5598.261

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H132").Value = 1807.6333
$ws.Range("I132").Value = 1667.16
$ws.Range("K132").Value = 5001.48
$ws.Range("M132").Value = -2471.48
$ws.Range("H135").Value = 39699
$ws.Range("J135").Value = 39699
$ws.Range("L135").Value = 39699
$ws.Range("N135").Value = -49839
$ws.Range("H136").Value = 7989.6562
$ws.Range("I136").Value = 6236.3335
$ws.Range("K136").Value = 18709.0005
$ws.Range("M136").Value = -16159.0005
